$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.488.46'
$ws.Range("E2").Value = '  +3.94%  '

$ws.Range("D3").Value = '1.803.11'
$ws.Range("E3").Value = '  +4.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3802'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3486'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.83'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.212'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07593'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.539'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.65%  '

$ws.Range("D15").Value = '1.800.90'
$ws.Range("E15").Value = '  +4.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.101'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06688'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.485'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.16%  '

$ws.Range("D23").Value = '27.445.51'
$ws.Range("E23").Value = '  +3.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.451'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.585'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.465'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.83%  '

$ws.Range("D30").Value = '2.005.51'
$ws.Range("E30").Value = '  +4.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.088'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.125'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08671'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.99%  '

$ws.Range("E36").Value = '  -1.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.517'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6882'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.73%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2224'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02370'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.46%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06392'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.57%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.877'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.278'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.07%  '

$ws.Range("E45").Value = '  +0.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6459'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.839'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.142'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07232'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.28%  '
